$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "06/08/2025"
$ws.Range("A23").Style = "Normal"
$ws.Range("B23").Value = "Binacional"
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = "Alianza Atl."
$ws.Range("F23").Value = "L"
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 2
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0.51
$ws.Range("L23").Value = 1.5
$ws.Range("M23").Value = 10
$ws.Range("N23").Value = 16
$ws.Range("O23").Value = 1
$ws.Range("P23").Value = 7
